$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so numeric-looking
# strings (e.g. "1.04", "91.151.72") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '91.151.72'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '3.179.91'
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '218.46'
$ws.Range('E5').Value = '  +2.36%  '
$ws.Range('D6').Value = '629.77'
$ws.Range('E6').Value = '  +1.00%  '
$ws.Range('D7').Value = '1.04'
$ws.Range('E7').Value = '  +25.23%  '
$ws.Range('D8').Value = '0.376'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('D9').Value = '0.999'
$ws.Range('D10').Value = '3.173.24'
$ws.Range('E10').Value = '  +2.42%  '
$ws.Range('D11').Value = '0.750'
$ws.Range('E11').Value = '  +22.65%  '
$ws.Range('E12').Value = '  +8.58%  '
$ws.Range('D13').Value = '0.0000252'
$ws.Range('E13').Value = '  +3.97%  '
$ws.Range('D14').Value = '35.27'
$ws.Range('E14').Value = '  +8.94%  '
$ws.Range('D15').Value = '5.57'
$ws.Range('E15').Value = '  +4.64%  '
$ws.Range('D16').Value = '91.005.27'
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').Value = '3.747.95'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').Value = '3.138.54'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').Value = '3.77'
$ws.Range('E19').Value = '  +11.24%  '
$ws.Range('D20').Value = '0.0000220'
$ws.Range('E20').Value = '  +2.02%  '
$ws.Range('D21').Value = '14.40'
$ws.Range('E21').Value = '  +6.74%  '
$ws.Range('D22').Value = '448.55'
$ws.Range('E22').Value = '  +5.06%  '
$ws.Range('D23').Value = '9.01'
$ws.Range('E23').Value = '  +8.59%  '
$ws.Range('D24').Value = '5.24'
$ws.Range('E24').Value = '  +5.68%  '
$ws.Range('D25').Value = '6.15'
$ws.Range('E25').Value = '  +11.17%  '
$ws.Range('D26').Value = '89.11'
$ws.Range('E26').Value = '  +5.91%  '
$ws.Range('D27').Value = '12.46'
$ws.Range('E27').Value = '  +3.38%  '
$ws.Range('D28').Value = '3.319.42'
$ws.Range('E28').Value = '  +1.82%  '
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').Value = '0.165'
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('D31').Value = '9.24'
$ws.Range('E31').Value = '  +13.29%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').Value = '0.982'
$ws.Range('E32').Value = '  -9.59%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '533.66'
$ws.Range('E33').Value = '  +4.43%  '
$ws.Range('D34').Value = '25.62'
$ws.Range('E34').Value = '  +13.96%  '
$ws.Range('D35').Value = '3.76'
$ws.Range('E35').Value = '  +1.52%  '
$ws.Range('D36').Value = '7.10'
$ws.Range('E36').Value = '  +5.72%  '
$ws.Range('D37').Value = '0.144'
$ws.Range('E37').Value = '  +10.92%  '
$ws.Range('D38').Value = '1.92'
$ws.Range('E38').Value = '  +6.32%  '
$ws.Range('D39').Value = '1.32'
$ws.Range('E39').Value = '  +4.68%  '
$ws.Range('E40').Value = '  -0.25%  '
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').Value = '0.156'
$ws.Range('E42').Value = '  +15.24%  '
$ws.Range('D43').Value = '0.409'
$ws.Range('E43').Value = '  +12.16%  '
$ws.Range('D44').Value = '0.0823'
$ws.Range('E44').Value = '  +14.52%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value = '1.96'
$ws.Range('E45').Value = '  +6.00%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').Value = '148.88'
$ws.Range('E47').Value = '  +2.03%  '
$ws.Range('D48').Value = '1.35'
$ws.Range('E48').Value = '  +9.38%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').Value = '4.45'
$ws.Range('E49').Value = '  +11.83%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '44.28'
$ws.Range('E50').Value = '  +1.82%  '
$ws.Range('D51').Value = '173.53'
$ws.Range('E51').Value = '  +8.00%  '
